$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F12").Value = 4955
$wsExhibition.Range("F17").Value = 352
$wsExhibition.Range("F20").Value = 284
$wsExhibition.Range("F22").Value = 3454

# Sheet "演出" (Performances)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 66

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 66
$wsAll.Range("F13").Value = 4955
$wsAll.Range("F18").Value = 352
$wsAll.Range("F21").Value = 284
$wsAll.Range("F23").Value = 3454
